$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the formatting of the
# existing header cells (e.g. G1 / "sum") by copying its formats over.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for the new column in row 2.
$ws.Range("H2").Value = 0
